$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$tbl = $ws.ListObjects.Item(1)

# ------------------------------------------------------------------
# 1) Update existing rows: H09, H11 now reference "Conductor"/"conductor"
#    instead of "propietario", and the H11 description text is updated
#    to talk about "el Conductor" instead of "el propietario".
# ------------------------------------------------------------------
$ws.Range("B10").Value = "Conductor"
$ws.Range("B12").Value = "conductor"
$ws.Range("E12").Value = "el Conductor esta dejando subir al bus a los usuarios, peor antes revisa en la lista de asientos vendidos  que las personas que estan ingresando si hayan pagado si no hacer el respectivo cobro, ademas de revisar que hayan hecho el chek in 0 contacto o si no proceder hacerlo antes de que suban al autobus"

# ------------------------------------------------------------------
# 2) Add a brand new row to the "RF" table for the new "Centro de
#    operaciones" story (H13) that attends PQR requests.
# ------------------------------------------------------------------
$newRow = $tbl.ListRows.Add()
$newRowIndex = $newRow.Range.Row

$ws.Cells.Item($newRowIndex, 1).Value = "H13"
$ws.Cells.Item($newRowIndex, 2).Value = "Centro de operaciones"
$ws.Cells.Item($newRowIndex, 3).Value = "atender PQR"
$ws.Cells.Item($newRowIndex, 4).Value = "poder dar respuesta a sus usuarios oconductores"
$ws.Cells.Item($newRowIndex, 5).Value = "el centro de operaciones reciben notificaciones de los usarios o conductores, este puede verlos y depues de esto puede darle una respuesta o generar una solucion con respecto a este PQR."

# Match the row height used by similar (H13 / 4-line) rows.
$ws.Rows.Item($newRowIndex).RowHeight = 57.6

# Copy cell formatting (borders/fill/number format) from the row above,
# which has the closest layout, then re-apply the left/center alignment
# used throughout columns B and C (see step 3).
$ws.Range("A" + ($newRowIndex - 1) + ":E" + ($newRowIndex - 1)).Copy() | Out-Null
$ws.Range("A" + $newRowIndex + ":E" + $newRowIndex).PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Re-write the values (PasteSpecial(Formats) does not touch values, but
# being explicit keeps this robust regardless of paste semantics).
$ws.Cells.Item($newRowIndex, 1).Value = "H13"
$ws.Cells.Item($newRowIndex, 2).Value = "Centro de operaciones"
$ws.Cells.Item($newRowIndex, 3).Value = "atender PQR"
$ws.Cells.Item($newRowIndex, 4).Value = "poder dar respuesta a sus usuarios oconductores"
$ws.Cells.Item($newRowIndex, 5).Value = "el centro de operaciones reciben notificaciones de los usarios o conductores, este puede verlos y depues de esto puede darle una respuesta o generar una solucion con respecto a este PQR."

# ------------------------------------------------------------------
# 3) Columns "YO COMO" (B) and "PUEDO" (C) are now left aligned and
#    vertically centered across the whole table (header included).
# ------------------------------------------------------------------
$lastRow = $tbl.Range.Rows.Count + $tbl.Range.Row - 1
$alignRange = $ws.Range("B1:C" + $lastRow)
$alignRange.HorizontalAlignment = -4131
$alignRange.VerticalAlignment = -4108

# ------------------------------------------------------------------
# 4) Update the selection to reflect the new end of the table.
# ------------------------------------------------------------------
$ws.Range("A17").Select()
